# Worklog Template.xlsx — add a week's worth of "Artefact Implementation"
# entries (4 Sep 2021 - 12 Sep 2021) continuing the pattern already present
# in rows 3-66, and move the sheet's on-screen selection down to where the
# new data ends (matches the author's saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$student     = "Phuong Dang"
$supervisor  = "Sasan Adibi"
$projectName = "P78-Data analytics and Deep Learning in better understanding of COVID-39"
$activity    = "Artefact Implementation"
$description = "Continue building the artefact, testing hyperparameters and fixing errors."

# Row 67 starts the day after the last filled-in row (66 -> 3 Sep 2021),
# and continues through row 75 (12 Sep 2021), mirroring rows 62-66 which
# already use the same project/activity/description combo.
$startRow  = 67
$endRow    = 75
$startDate = $ws.Cells.Item(66, 4).Value2 + 1

for ($i = 0; $i -lt ($endRow - $startRow + 1); $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $student
    $ws.Cells.Item($r, 2).Value2 = $supervisor
    $ws.Cells.Item($r, 3).Value2 = $projectName
    $ws.Cells.Item($r, 4).Value2 = $startDate + $i
    $ws.Cells.Item($r, 5).Value2 = 120
    $ws.Cells.Item($r, 6).Value2 = $activity
    $ws.Cells.Item($r, 7).Value2 = $description
}

# Move the selection to where the author left off editing.
$ws.Range("G76").Select()
